$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

$ws.Range("A2").Value = "Employee 1"
$ws.Range("B2").Value = "OFF"
$ws.Range("A3").Value = "Employee 2"
$ws.Range("A4").Value = "Employee 3"
$ws.Range("A5").Value = "Employee 4"
$ws.Range("A6").Value = "Employee 5"
$ws.Range("A7").Value = "Employee 6"

$ws.Range("D9").Select()
